# Atualizado por script em 11-11-2023 20:45
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: swap the F:V content of two rows (keeps A:E columns,
# i.e. Indice/pais/torneio/temporada/data_partida, untouched).
function Swap-RowData($rowA, $rowB) {
    $dataA = $ws.Range("F" + $rowA + ":V" + $rowA).Value2
    $dataB = $ws.Range("F" + $rowB + ":V" + $rowB).Value2
    $ws.Range("F" + $rowA + ":V" + $rowA).Value2 = $dataB
    $ws.Range("F" + $rowB + ":V" + $rowB).Value2 = $dataA
}

# Rows 89 / 90 swap (Trabzonspor-Alanyaspor <-> F.Karagumruk-Adana Demirspor)
Swap-RowData 89 90

# Rows 91 / 92 swap (Hatayspor-Kayserispor <-> Kasimpasa-Istanbulspor AS)
Swap-RowData 91 92

# Rows 94 / 95 swap (Adana Demirspor-Konyaspor <-> Rizespor-Galatasaray)
Swap-RowData 94 95

# Rows 97 / 98 / 99 rotate: new97 = old99, new98 = old97, new99 = old98
$data97 = $ws.Range("F97:V97").Value2
$data98 = $ws.Range("F98:V98").Value2
$data99 = $ws.Range("F99:V99").Value2
$ws.Range("F97:V97").Value2 = $data99
$ws.Range("F98:V98").Value2 = $data97
$ws.Range("F99:V99").Value2 = $data98

# New row 117: Hatayspor 2 - 1 Galatasaray
# Copy formatting (styles/number formats) from the last existing row (116).
$ws.Range("A116:V116").Copy($ws.Range("A117:V117"))

$ws.Cells.Item(117, 1).Value2 = 116
$ws.Cells.Item(117, 2).Value2 = "turkey"
$ws.Cells.Item(117, 3).Value2 = "super-lig"
$ws.Cells.Item(117, 4).Value2 = "2023-2024"
$ws.Cells.Item(117, 5).Value2 = 45241.70833333334
$ws.Cells.Item(117, 6).Value2 = "Hatayspor"
$ws.Cells.Item(117, 7).Value2 = 2
$ws.Cells.Item(117, 8).Value2 = "Galatasaray"
$ws.Cells.Item(117, 9).Value2 = 1
$ws.Cells.Item(117, 10).Value2 = 6.26
$ws.Cells.Item(117, 11).Value2 = "04/11/2023 17:13"
$ws.Cells.Item(117, 12).Value2 = 6.54
$ws.Cells.Item(117, 13).Value2 = "11/11/2023 16:57"
$ws.Cells.Item(117, 14).Value2 = 5.09
$ws.Cells.Item(117, 15).Value2 = "04/11/2023 17:13"
$ws.Cells.Item(117, 16).Value2 = 5.14
$ws.Cells.Item(117, 17).Value2 = "11/11/2023 16:57"
$ws.Cells.Item(117, 18).Value2 = 1.47
$ws.Cells.Item(117, 19).Value2 = "04/11/2023 17:13"
$ws.Cells.Item(117, 20).Value2 = 1.46
$ws.Cells.Item(117, 21).Value2 = "11/11/2023 16:57"
$ws.Cells.Item(117, 22).Value2 = "https://www.betexplorer.com/football/turkey/super-lig/hatayspor-galatasaray/8AaXwqs6/"
